# "updated schedule and Android App"
# The burndown task list is reworked: the old coarse-grained tasks
# ("Hardware Prototyp" / "Software App") are replaced by a more detailed
# breakdown of hardware and software sub-tasks, and the schedule rows are
# reordered/renumbered accordingly. Row 2 is intentionally left blank in
# the new layout (data now starts at row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the previous task rows (everything below the header) before rewriting.
$ws.Range("A2:F5").ClearContents()

$tasks = @(
    @{ Row = 3;  Name = "Harware Selection";               Cur = 3;  Org = 3;  Eff = 0; Rem = 3 },
    @{ Row = 4;  Name = "Hardware Assembly";                Cur = 2;  Org = 2;  Eff = 0; Rem = 2 },
    @{ Row = 5;  Name = "Hardware Enclosure";                Cur = 5;  Org = 5;  Eff = 0; Rem = 5 },
    @{ Row = 6;  Name = "Software Technology Selection";     Cur = 5;  Org = 3;  Eff = 3; Rem = 0 },
    @{ Row = 7;  Name = "Software ""Hello World""";          Cur = 5;  Org = 5;  Eff = 0; Rem = 5 },
    @{ Row = 8;  Name = "Software Coding";                   Cur = 30; Org = 30; Eff = 0; Rem = 30 },
    @{ Row = 9;  Name = "Firmware Prototyp";                 Cur = 20; Org = 20; Eff = 0; Rem = 20 },
    @{ Row = 10; Name = "Software Testing";                  Cur = 30; Org = 30; Eff = 0; Rem = 30 },
    @{ Row = 11; Name = "Testing/Shredding";                 Cur = 30; Org = 30; Eff = 0; Rem = 30 }
)

foreach ($t in $tasks) {
    $r = $t.Row
    $ws.Cells.Item($r, 1).Value = $t.Name   # Task
    $ws.Cells.Item($r, 2).Value = $t.Cur    # Cur. Est.
    $ws.Cells.Item($r, 3).Value = $t.Org    # Org. Est.
    $ws.Cells.Item($r, 4).Value = $t.Eff    # Effort
    $ws.Cells.Item($r, 5).Value = $t.Rem    # Remain
}

# Column A now needs to fit the longer task names.
$ws.Columns.Item(1).AutoFit()

# Mirror the saved selection state from the edited file.
$ws.Range("F6").Select()
